# Update the division problems in the practice-sheet table.
# Each data row of the table is followed by three blank spacer rows,
# so the populated rows are 1, 5, 9, 13, 17 (1-indexed).  We address
# cells directly by (row, column) and replace the cell's Range.Text
# outright (rather than Range.Find, whose scoping isn't reliably
# confined to a sub-range here) to avoid any ambiguity from values
# (e.g. "94÷9=") that repeat elsewhere in the sheet but map to
# different replacements.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$edits = @(
    @{ Row = 1;  Col = 2; New = "42÷5=" },
    @{ Row = 1;  Col = 3; New = "75÷8=" },
    @{ Row = 1;  Col = 4; New = "12÷8=" },
    @{ Row = 1;  Col = 5; New = "35÷3=" },

    @{ Row = 5;  Col = 1; New = "21÷5=" },
    @{ Row = 5;  Col = 2; New = "81÷9=" },
    @{ Row = 5;  Col = 3; New = "29÷7=" },
    @{ Row = 5;  Col = 4; New = "11÷4=" },
    @{ Row = 5;  Col = 5; New = "56÷9=" },

    @{ Row = 9;  Col = 1; New = "51÷9=" },
    @{ Row = 9;  Col = 2; New = "27÷9=" },
    @{ Row = 9;  Col = 3; New = "91÷2=" },
    @{ Row = 9;  Col = 4; New = "83÷7=" },
    @{ Row = 9;  Col = 5; New = "85÷8=" },

    @{ Row = 13; Col = 1; New = "32÷6=" },
    @{ Row = 13; Col = 2; New = "67÷2=" },
    @{ Row = 13; Col = 3; New = "97÷8=" },
    @{ Row = 13; Col = 4; New = "68÷8=" },
    @{ Row = 13; Col = 5; New = "82÷5=" },

    @{ Row = 17; Col = 1; New = "40÷2=" },
    @{ Row = 17; Col = 2; New = "66÷6=" },
    @{ Row = 17; Col = 3; New = "22÷7=" },
    @{ Row = 17; Col = 4; New = "72÷2=" },
    @{ Row = 17; Col = 5; New = "26÷8=" }
)

foreach ($edit in $edits) {
    $cell = $t.Cell($edit.Row, $edit.Col)
    $rng = $cell.Range
    # Drop the trailing cell-end mark so only the visible text is replaced.
    $rng.MoveEnd(12, -1) | Out-Null
    $rng.Text = $edit.New
}
